$wb = $excel.ActiveWorkbook

# --- Worksheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 1423.0769
$ws.Range("I9").Value = 1137.5
$ws.Range("J9").Value = 1880
$ws.Range("K9").Value = 1137.5
$ws.Range("L9").Value = 1880
$ws.Range("M9").Value = -968.5
$ws.Range("N9").Value = -2218
$ws.Range("H80").Value = 45722.363
$ws.Range("I80").Value = 20219.4
$ws.Range("J80").Value = 66974.836
$ws.Range("K80").Value = 60658.2
$ws.Range("L80").Value = 200924.508
$ws.Range("M80").Value = -59660.2
$ws.Range("N80").Value = -202920.508
$ws.Range("H83").Value = 45722.363
$ws.Range("I83").Value = 20219.4
$ws.Range("J83").Value = 66974.836
$ws.Range("K83").Value = 181974.6
$ws.Range("L83").Value = 602773.524
$ws.Range("M83").Value = -176982.6
$ws.Range("N83").Value = -612757.524
$ws.Range("H100").Value = 4757.3335
$ws.Range("I100").Value = 4671.727
$ws.Range("J100").Value = 5699
$ws.Range("K100").Value = 4671.727
$ws.Range("L100").Value = 5699
$ws.Range("M100").Value = -4130.727
$ws.Range("N100").Value = -6781
$ws.Range("H106").Value = 2459.0715
$ws.Range("I106").Value = 2306.3845
$ws.Range("J106").Value = 4444
$ws.Range("K106").Value = 2306.3845
$ws.Range("L106").Value = 4444
$ws.Range("M106").Value = -1675.3845
$ws.Range("N106").Value = -5706
$ws.Range("H113").Value = 134261940
$ws.Range("I113").Value = 138891060
$ws.Range("J113").Value = 125003690
$ws.Range("K113").Value = 138891060
$ws.Range("L113").Value = 125003690
$ws.Range("M113").Value = -138887806
$ws.Range("N113").Value = -125010198
$ws.Range("H132").Value = 1894.0702
$ws.Range("I132").Value = 1561.6123
$ws.Range("J132").Value = 3930.375
$ws.Range("K132").Value = 4684.8369
$ws.Range("L132").Value = 11791.125
$ws.Range("M132").Value = -2154.8369
$ws.Range("N132").Value = -16851.125
$ws.Range("H138").Value = 1520367
$ws.Range("I138").Value = 3128.3333
$ws.Range("J138").Value = 2089331.4
$ws.Range("K138").Value = 9384.999899999999
$ws.Range("L138").Value = 6267994.199999999
$ws.Range("M138").Value = -4244.999899999999
$ws.Range("N138").Value = -6278274.199999999

# --- Worksheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3009.3457
$ws.Range("I32").Value = 2957.3948
$ws.Range("J32").Value = 3799
$ws.Range("K32").Value = 2957.3948
$ws.Range("L32").Value = 3799
$ws.Range("M32").Value = -2670.3948
$ws.Range("N32").Value = -4373
$ws.Range("H43").Value = 5999.5
$ws.Range("I43").Value = 9999
$ws.Range("J43").Value = 2000
$ws.Range("K43").Value = 9999
$ws.Range("L43").Value = 2000
$ws.Range("M43").Value = -9686
$ws.Range("N43").Value = -2626
$ws.Range("H74").Value = 43310.95
$ws.Range("I74").Value = 53248.387
$ws.Range("J74").Value = 4803.375
$ws.Range("K74").Value = 53248.387
$ws.Range("L74").Value = 4803.375
$ws.Range("M74").Value = -52374.387
$ws.Range("N74").Value = -6551.375
$ws.Range("H77").Value = 43310.95
$ws.Range("I77").Value = 53248.387
$ws.Range("J77").Value = 4803.375
$ws.Range("K77").Value = 266241.935
$ws.Range("L77").Value = 24016.875
$ws.Range("M77").Value = -261873.935
$ws.Range("N77").Value = -32752.875
$ws.Range("H109").Value = 58880
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 58880
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 58880
$ws.Range("N109").Value = -61654
$ws.Range("H132").Value = 733978.9399999999
$ws.Range("I132").Value = 1120985.5
$ws.Range("J132").Value = 6406.44
$ws.Range("K132").Value = 3362956.5
$ws.Range("L132").Value = 19219.32
$ws.Range("M132").Value = -3360426.5
$ws.Range("N132").Value = -24279.32

# --- Worksheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3273.0908
$ws.Range("I105").Value = 1361.625
$ws.Range("J105").Value = 8370.333000000001
$ws.Range("K105").Value = 1361.625
$ws.Range("L105").Value = 8370.333000000001
$ws.Range("M105").Value = 385.375
$ws.Range("N105").Value = -11864.333

# --- Worksheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 1125.25
$ws.Range("I19").Value = 1125.25
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 1125.25
$ws.Range("L19").Value = 0
$ws.Range("M19").ClearContents()
$ws.Range("N19").Value = -955.25
$ws.Range("H24").Value = 1125.25
$ws.Range("I24").Value = 1125.25
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 1125.25
$ws.Range("L24").Value = 0
$ws.Range("M24").ClearContents()
$ws.Range("N24").Value = -955.25
$ws.Range("H58").Value = 9096239
$ws.Range("I58").Value = 20001856
$ws.Range("J58").Value = 8224.799999999999
$ws.Range("K58").Value = 20001856
$ws.Range("L58").Value = 8224.799999999999
$ws.Range("M58").Value = -20001653
$ws.Range("N58").Value = -8630.799999999999
$ws.Range("H132").Value = 3210.8406
$ws.Range("I132").Value = 1756.3462
$ws.Range("J132").Value = 7659.8823
$ws.Range("K132").Value = 5269.0386
$ws.Range("L132").Value = 22979.6469
$ws.Range("M132").Value = -2739.0386
$ws.Range("N132").Value = -28039.6469
$ws.Range("H134").Value = 4339.5615
$ws.Range("I134").Value = 1765.15
$ws.Range("J134").Value = 10397
$ws.Range("K134").Value = 5295.450000000001
$ws.Range("L134").Value = 31191
$ws.Range("M134").Value = -2760.450000000001
$ws.Range("N134").Value = -36261
$ws.Range("H136").Value = 9096239
$ws.Range("I136").Value = 20001856
$ws.Range("J136").Value = 8224.799999999999
$ws.Range("K136").Value = 60005568
$ws.Range("L136").Value = 24674.4
$ws.Range("M136").Value = -60003018
$ws.Range("N136").Value = -29774.4

# --- Worksheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 80118310
$ws.Range("I4").Value = 71958840
$ws.Range("J4").Value = 92728410
$ws.Range("K4").Value = 215876520
$ws.Range("L4").Value = 278185230
$ws.Range("M4").Value = -215876408
$ws.Range("N4").Value = -278185454
$ws.Range("H34").Value = 6751.643
$ws.Range("I34").Value = 88.333336
$ws.Range("J34").Value = 11749.125
$ws.Range("K34").Value = 265.000008
$ws.Range("L34").Value = 35247.375
$ws.Range("M34").Value = -181.000008
$ws.Range("N34").Value = -35415.375
$ws.Range("H68").Value = 3701.8076
$ws.Range("I68").Value = 2374.8333
$ws.Range("J68").Value = 4099.9
$ws.Range("K68").Value = 7124.499899999999
$ws.Range("L68").Value = 12299.7
$ws.Range("M68").Value = -6313.499899999999
$ws.Range("N68").Value = -13921.7
$ws.Range("H71").Value = 3701.8076
$ws.Range("I71").Value = 2374.8333
$ws.Range("J71").Value = 4099.9
$ws.Range("K71").Value = 21373.4997
$ws.Range("L71").Value = 36899.1
$ws.Range("M71").Value = -17317.4997
$ws.Range("N71").Value = -45011.1
$ws.Range("H104").Value = 4944.6665
$ws.Range("I104").Value = 0
$ws.Range("J104").Value = 4944.6665
$ws.Range("K104").Value = 0
$ws.Range("L104").ClearContents()
$ws.Range("M104").Value = 14833.9995
$ws.Range("N104").Value = -20075.9995
$ws.Range("H121").Value = 1440
$ws.Range("I121").Value = 480.42856
$ws.Range("J121").Value = 4798.5
$ws.Range("K121").Value = 1441.28568
$ws.Range("L121").Value = 14395.5
$ws.Range("M121").Value = -131.28568
$ws.Range("N121").Value = -17015.5
$ws.Range("H131").Value = 1463.8334
$ws.Range("I131").Value = 976.6
$ws.Range("J131").Value = 3900
$ws.Range("K131").Value = 2929.8
$ws.Range("L131").Value = 11700
$ws.Range("M131").Value = 2110.2
$ws.Range("N131").Value = -21780
$ws.Range("H136").Value = 2512.5
$ws.Range("I136").Value = 2016.6666
$ws.Range("J136").Value = 4000
$ws.Range("K136").Value = 6049.9998
$ws.Range("L136").Value = 12000
$ws.Range("M136").Value = -949.9997999999996
$ws.Range("N136").Value = -22200

# --- Worksheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 10600.308
$ws.Range("I70").Value = 9472.833000000001
$ws.Range("J70").Value = 11566.714
$ws.Range("K70").Value = 9472.833000000001
$ws.Range("L70").Value = 11566.714
$ws.Range("M70").Value = -9202.833000000001
$ws.Range("N70").Value = -12106.714
$ws.Range("H73").Value = 10600.308
$ws.Range("I73").Value = 9472.833000000001
$ws.Range("J73").Value = 11566.714
$ws.Range("K73").Value = 9472.833000000001
$ws.Range("L73").Value = 11566.714
$ws.Range("M73").Value = -8536.833000000001
$ws.Range("N73").Value = -13438.714
$ws.Range("H97").Value = 1063.5
$ws.Range("I97").Value = 1070.6552
$ws.Range("J97").Value = 1022
$ws.Range("K97").Value = 1070.6552
$ws.Range("L97").Value = 1022
$ws.Range("M97").Value = -574.6551999999999
$ws.Range("N97").Value = -2014
$ws.Range("H122").Value = 55615176
$ws.Range("I122").Value = 91002770
$ws.Range("J122").Value = 6109.143
$ws.Range("K122").Value = 273008310
$ws.Range("L122").Value = 18327.429
$ws.Range("M122").Value = -273005860
$ws.Range("N122").Value = -23227.429

# --- Worksheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5967.15
$ws.Range("I40").Value = 4611.5835
$ws.Range("J40").Value = 8000.5
$ws.Range("K40").Value = 4611.5835
$ws.Range("L40").Value = 8000.5
$ws.Range("M40").Value = -4475.5835
$ws.Range("N40").Value = -8272.5
$ws.Range("H122").Value = 8126.8
$ws.Range("I122").Value = 7179.6
$ws.Range("J122").Value = 8600.4
$ws.Range("K122").Value = 21538.8
$ws.Range("L122").Value = 25801.2
$ws.Range("M122").Value = -19088.8
$ws.Range("N122").Value = -30701.2
$ws.Range("H132").Value = 20008970
$ws.Range("I132").Value = 50003948
$ws.Range("J132").Value = 12317
$ws.Range("K132").Value = 150011844
$ws.Range("L132").Value = 36951
$ws.Range("M132").Value = -150009314
$ws.Range("N132").Value = -42011
$ws.Range("H136").Value = 7717.25
$ws.Range("I136").Value = 3014.7
$ws.Range("J136").Value = 19473.625
$ws.Range("K136").Value = 9044.099999999999
$ws.Range("L136").Value = 58420.875
$ws.Range("M136").Value = -6494.099999999999
$ws.Range("N136").Value = -63520.875

# --- Worksheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 47620548
$ws.Range("I107").Value = 1800.2
$ws.Range("J107").Value = 166667420
$ws.Range("K107").Value = 5400.6
$ws.Range("L107").Value = 500002260
$ws.Range("M107").Value = -3480.6
$ws.Range("N107").Value = -500006100
$ws.Range("H113").Value = 684.8333
$ws.Range("I113").Value = 284
$ws.Range("J113").Value = 1085.6666
$ws.Range("K113").Value = 852
$ws.Range("L113").Value = 3256.9998
$ws.Range("M113").Value = 1318
$ws.Range("N113").Value = -7596.9998
$ws.Range("H132").Value = 7366722
$ws.Range("I132").Value = 9097986
$ws.Range("J132").Value = 42143.31
$ws.Range("K132").Value = 27293958
$ws.Range("L132").Value = 126429.93
$ws.Range("M132").Value = -27291428
$ws.Range("N132").Value = -131489.93
